# Update the cryptocurrency price/volume table (columns D and E, rows 2-51)
# to reflect the latest scrape, matching the committed OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.338.38'
$ws.Range('E2').Value = '  -4.06%  '
$ws.Range('D3').Value = '1.860.21'
$ws.Range('E3').Value = '  -5.01%  '
$ws.Range('E4').Value = '  -0.95%  '
$ws.Range('D5').Value = '''323.42'
$ws.Range('E5').Value = '  +0.34%  '
$ws.Range('E6').Value = '  -0.93%  '
$ws.Range('E7').Value = '  -5.82%  '
$ws.Range('D8').Value = '''0.3861'
$ws.Range('E8').Value = '  -5.19%  '
$ws.Range('D9').Value = '''48.00'
$ws.Range('E9').Value = '  -10.99%  '
$ws.Range('D10').Value = '''0.07897'
$ws.Range('E10').Value = '  -7.12%  '
$ws.Range('D11').Value = '''1.021'
$ws.Range('E11').Value = '  -3.65%  '
$ws.Range('E12').Value = '  -4.80%  '
$ws.Range('D13').Value = '1.858.76'
$ws.Range('E13').Value = '  -5.39%  '
$ws.Range('D14').Value = '''5.892'
$ws.Range('E14').Value = '  -4.49%  '
$ws.Range('D15').Value = '''7.143'
$ws.Range('E15').Value = '  -5.81%  '
$ws.Range('D16').Value = '''1.000'
$ws.Range('E16').Value = '  -1.21%  '
$ws.Range('D17').Value = '''0.00001035'
$ws.Range('E17').Value = '  -3.61%  '
$ws.Range('D18').Value = '''85.61'
$ws.Range('E18').Value = '  -5.55%  '
$ws.Range('D19').Value = '''0.06531'
$ws.Range('E19').Value = '  -1.52%  '
$ws.Range('D20').Value = '''17.18'
$ws.Range('E20').Value = '  -7.15%  '
$ws.Range('D21').Value = '''1.001'
$ws.Range('E21').Value = '  -0.92%  '
$ws.Range('E22').Value = '  -5.53%  '
$ws.Range('D23').Value = '27.343.39'
$ws.Range('E23').Value = '  -4.10%  '
$ws.Range('E24').Value = '  -5.05%  '
$ws.Range('E25').Value = '  -0.89%  '
$ws.Range('D26').Value = '2.081.90'
$ws.Range('E26').Value = '  -5.26%  '
$ws.Range('D27').Value = '''152.76'
$ws.Range('E27').Value = '  -2.38%  '
$ws.Range('D28').Value = '''19.79'
$ws.Range('E28').Value = '  -2.69%  '
$ws.Range('D29').Value = '''2.061'
$ws.Range('E29').Value = '  -5.23%  '
$ws.Range('D30').Value = '''5.486'
$ws.Range('E30').Value = '  -6.06%  '
$ws.Range('D31').Value = '''120.61'
$ws.Range('E31').Value = '  -3.06%  '
$ws.Range('D32').Value = '''1.490'
$ws.Range('E32').Value = '  +2.33%  '
$ws.Range('D33').Value = '''0.09305'
$ws.Range('E33').Value = '  -3.69%  '
$ws.Range('D34').Value = '''0.9313'
$ws.Range('E34').Value = '  -5.51%  '
$ws.Range('D35').Value = '''3.601'
$ws.Range('E35').Value = '  -2.45%  '
$ws.Range('D36').Value = '''5.266'
$ws.Range('E36').Value = '  -6.66%  '
$ws.Range('D37').Value = '''0.02231'
$ws.Range('E37').Value = '  -4.48%  '
$ws.Range('D38').Value = '''0.05996'
$ws.Range('E38').Value = '  -3.21%  '
$ws.Range('D39').Value = '''1.220'
$ws.Range('E39').Value = '  -2.67%  '
$ws.Range('D40').Value = '''8.254'
$ws.Range('E40').Value = '  -9.42%  '
$ws.Range('E41').Value = '  -1.01%  '
$ws.Range('D42').Value = '''0.5909'
$ws.Range('E42').Value = '  -5.28%  '
$ws.Range('D43').Value = '''0.1887'
$ws.Range('E43').Value = '  -1.60%  '
$ws.Range('D44').Value = '''10.14'
$ws.Range('E44').Value = '  -9.51%  '
$ws.Range('D45').Value = '''1.275'
$ws.Range('E45').Value = '  -5.83%  '
$ws.Range('D46').Value = '''0.5619'
$ws.Range('E46').Value = '  -5.69%  '
$ws.Range('D47').Value = '''12.03'
$ws.Range('E47').Value = '  -7.42%  '
$ws.Range('E48').Value = '  -1.25%  '
$ws.Range('D49').Value = '''1.925'
$ws.Range('E49').Value = '  -6.71%  '
$ws.Range('D50').Value = '''0.06785'
$ws.Range('E50').Value = '  -0.46%  '
$ws.Range('D51').Value = '''107.82'
$ws.Range('E51').Value = '  -3.07%  '
